$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Pet" section header text (now spans a wider area, O1:R1)
$ws.Range("O1").Value = "                                      Pet                       "

# Extend the header's formatting (fill/border) from Q1 into the new R1 cell
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)

# New "Nome" sub-header for the Pet table in R2, matching Q2's formatting
$ws.Range("Q2").Copy()
$ws.Range("R2").PasteSpecial(-4122)
$ws.Range("R2").Value = "Nome"

# New data rows for the Pet table column R
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)
$ws.Range("R3").Value = "Cachorro"

$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = "Gato"

# Update the selection to match the saved view state (column N selected)
$ws.Columns("N").Select() | Out-Null
